# The source dataset gained one new weekly observation at the top of the
# Pina (Caramelo / Tercera) block for "Feria Lagunitas de Puerto Montt".
# All existing records (rows 203-249) shift down by one row, and a brand
# new row is written in at row 203 with the newest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 203:249 down to 204:250, carrying values + formatting along
# (this also grows the sheet dimension to A1:T250 automatically).
$ws.Rows.Item(203).Insert()

# Populate the newly inserted row 203 with the new observation.
$ws.Range("A203").Value = 4
$ws.Range("B203").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C203").Value = 'Los Lagos'
$ws.Range("D203").Value = 44711
$ws.Range("E203").Value = 10
$ws.Range("F203").Value = 'Fruta'
$ws.Range("G203").Value = 100108
$ws.Range("H203").Value = 'Tropicales y subtropicales'
$ws.Range("I203").Value = 100108005
$ws.Range("J203").Value = 'Piña'
$ws.Range("K203").Value = 'Caramelo'
$ws.Range("L203").Value = 'Tercera'
$ws.Range("M203").Value = 80
$ws.Range("N203").Value = 21000
$ws.Range("O203").Value = 21000
$ws.Range("P203").Value = 21000
$ws.Range("Q203").Value = '$/caja 16 unidades'
$ws.Range("R203").Value = 'Ecuador'
$ws.Range("S203").Value = 1312
$ws.Range("T203").Value = 16
